# "Add files via upload" - re-uploaded workbook with a couple of small
# housekeeping fixes:
#   1. Drop the stray " (2)" suffix that had crept into the
#      "EtOAc Ethanol Water" tab name.
#   2. Leave the workbook with the "D-Limonene Ethanol Water" tab active
#      (instead of "EtOAc Ethanol Water (2)"), and resize/reposition the
#      window to the author's last-used layout.

$wb = $excel.ActiveWorkbook

# 1. Rename "EtOAc Ethanol Water (2)" -> "EtOAc Ethanol Water"
$etoac = $wb.Worksheets.Item("EtOAc Ethanol Water (2)")
$etoac.Name = "EtOAc Ethanol Water"

# 2. Make "D-Limonene Ethanol Water" the active/selected sheet.
$dlim = $wb.Worksheets.Item("D-Limonene Ethanol Water")
$dlim.Activate()
$dlim.Select()

# 3. Restore the saved window position/size for this view.
$aw = $excel.ActiveWindow
$aw.Left = 10104
$aw.Top = 420
$aw.Width = 12828
$aw.Height = 12624

# Scroll the sheet tab strip so the active sheet's group is in view.
$aw.ScrollWorkbookTabs(12, 0)
